$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new MARKUP row (row 3) below the existing settings
$ws.Range("A3").Value = "MARKUP"
$ws.Range("B3").Value = 0

# Move/keep the active selection where the author last left it
[void]$ws.Range("C6").Select()
